$wb = $excel.ActiveWorkbook

# --- Sprint 4 sheet: fill in Day columns (S, S, M) with the Start hours values ---
$sprint4 = $wb.Worksheets.Item("Sprint 4")

for ($row = 2; $row -le 6; $row++) {
    $startHours = $sprint4.Cells.Item($row, 3).Value2   # column C = "Start hours"
    if ($startHours -ne $null) {
        $sprint4.Cells.Item($row, 4).Value = $startHours  # D = S
        $sprint4.Cells.Item($row, 5).Value = $startHours  # E = S
        $sprint4.Cells.Item($row, 6).Value = $startHours  # F = M
    }
}

# Move selection on Sprint 4 and make it the visible/selected sheet
[void]$sprint4.Select()
[void]$sprint4.Range("G4").Select()

# --- Sprint 3 sheet: scroll the view down (no longer the tab-selected sheet) ---
$sprint3 = $wb.Worksheets.Item("Sprint 3")
[void]$sprint3.Activate()
$topLeft = $sprint3.Range("A13")
$excel.ActiveWindow.ScrollRow = $topLeft.Row
$excel.ActiveWindow.ScrollColumn = $topLeft.Column
[void]$sprint3.Range("B7").Select()

# Re-select Sprint 4 last so it ends up as the active/visible tab
[void]$sprint4.Select()
[void]$sprint4.Range("G4").Select()
